# Updated cryptos list with GitHub Actions data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $range = $ws.Range($cellRef)
    $range.Value = "'" + $text
    $range.ClearFormats()
}

Set-TextValue "D2" "27.961.91"
Set-TextValue "E2" "  -4.22%  "
Set-TextValue "D3" "1.740.61"
Set-TextValue "E3" "  -4.58%  "
Set-TextValue "D4" "1.001"
Set-TextValue "D5" "226.70"
Set-TextValue "D6" "0.5797"
Set-TextValue "E6" "  -3.25%  "
Set-TextValue "D7" "1.002"
Set-TextValue "E7" "  -0.10%  "
Set-TextValue "D8" "0.2739"
Set-TextValue "E8" "  -0.89%  "
Set-TextValue "E9" "  -1.24%  "
Set-TextValue "D10" "0.06631"
Set-TextValue "E10" "  -4.52%  "
Set-TextValue "D11" "0.07554"
Set-TextValue "E11" "  -0.59%  "
Set-TextValue "D12" "1.743.42"
Set-TextValue "E12" "  -4.63%  "
Set-TextValue "D13" "4.711"
Set-TextValue "E13" "  -0.28%  "
Set-TextValue "D14" "0.6024"
Set-TextValue "E14" "  -4.06%  "
Set-TextValue "D15" "1.977.47"
Set-TextValue "E15" "  -4.57%  "
Set-TextValue "D16" "74.70"
Set-TextValue "D17" "0.000008752"
Set-TextValue "E17" "  -10.83%  "
Set-TextValue "D18" "27.940.55"
Set-TextValue "E18" "  -3.66%  "
Set-TextValue "D19" "5.317"
Set-TextValue "E19" "  -3.92%  "
Set-TextValue "E20" "  -0.14%  "
Set-TextValue "D21" "205.62"
Set-TextValue "E21" "  -4.90%  "
Set-TextValue "D22" "11.29"
Set-TextValue "E22" "  -2.27%  "
Set-TextValue "D23" "6.632"
Set-TextValue "E23" "  -3.07%  "
Set-TextValue "E24" "  -0.20%  "
Set-TextValue "D25" "150.34"
Set-TextValue "E25" "  -3.39%  "
Set-TextValue "D26" "8.044"
Set-TextValue "E26" "  +1.25%  "
Set-TextValue "E27" "  -4.33%  "
Set-TextValue "E28" "  -1.83%  "
Set-TextValue "D29" "1.386"
Set-TextValue "E29" "  -2.67%  "
Set-TextValue "D30" "0.06180"
Set-TextValue "E30" "  -4.41%  "
Set-TextValue "E31" "  -3.24%  "
Set-TextValue "D32" "3.747"
Set-TextValue "E32" "  -1.77%  "
Set-TextValue "E33" "  -1.03%  "
Set-TextValue "D34" "1.676"
Set-TextValue "E34" "  -2.50%  "
Set-TextValue "E35" "  -5.01%  "
Set-TextValue "D36" "0.6404"
Set-TextValue "E36" "  -0.76%  "
Set-TextValue "E37" "  -4.89%  "
Set-TextValue "D38" "2.716"
Set-TextValue "E38" "  -1.10%  "
Set-TextValue "D39" "0.01672"
Set-TextValue "E39" "  -4.34%  "
Set-TextValue "D40" "1.124.26"
Set-TextValue "E40" "  -0.66%  "
Set-TextValue "D41" "6.154"
Set-TextValue "E41" "  -6.54%  "
Set-TextValue "D42" "0.8758"
Set-TextValue "E42" "  -1.67%  "
Set-TextValue "D43" "1.003"
Set-TextValue "E43" "  +0.08%  "
Set-TextValue "D44" "99.98"
Set-TextValue "E44" "  -0.55%  "
Set-TextValue "D45" "1.889.46"
Set-TextValue "E45" "  -4.74%  "
Set-TextValue "D46" "59.38"
Set-TextValue "E46" "  -4.47%  "
Set-TextValue "B47" "RenderToken"
Set-TextValue "C47" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D47" "1.580"
Set-TextValue "E47" "  -1.98%  "
Set-TextValue "B48" "BabyDogeCoin"
Set-TextValue "C48" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D48" "0.00000000107"
Set-TextValue "E48" "  -5.23%  "
Set-TextValue "D49" "8.259"
Set-TextValue "E49" "  -2.37%  "
Set-TextValue "E50" "  -2.26%  "
Set-TextValue "D51" "6.259"
Set-TextValue "E51" "  -1.53%  "
